$d = $word.ActiveDocument

# Locate the "課題６−４" heading, then the empty paragraph right after it —
# that is the paragraph that gains the pasted C source + program output.
$r = $d.Content
$found = $r.Find.Execute("課題６−４", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find heading 課題６−４"
}

$headingPara = $r.Paragraphs(1)
$target = $headingPara.Next()

$newXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="808080"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>#include</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="A31515"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>&lt;stdio.h&gt;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>int</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> main(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>void</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>){</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>int</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> i,a;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>int</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> wa = </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="098658"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>0</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">  </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>int</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> seki = </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="098658"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">  printf(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="A31515"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>"please enter an  integer: "</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>);</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">  scanf(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="A31515"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>"%d"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>, &amp;a);</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>for</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> (i=</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="098658"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>1</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;i&lt;=a;i++){</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">          wa = wa + i;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">          seki = seki * i;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">    }</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">        printf(</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="A31515"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>"Sum: %d\n"</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>, seki);</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve">        </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="0000FF"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>return</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="098658"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>0</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>;</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="0" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
        <w:t>}</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
        <w:spacing w:after="240" w:line="270" w:lineRule="atLeast"/>
        <w:rPr>
          <w:rFonts w:ascii="Menlo" w:eastAsia="Times New Roman" w:hAnsi="Menlo" w:cs="Menlo"/>
          <w:color w:val="000000"/>
          <w:kern w:val="0"/>
          <w:sz w:val="18"/>
          <w:szCs w:val="18"/>
          <w14:ligatures w14:val="none"/>
        </w:rPr>
      </w:pPr>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>please enter an  integer: 4</w:t>
      </w:r>
    </w:p><w:p>
      <w:pPr>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:lang w:val="en-US"/>
        </w:rPr>
        <w:t>Sum: 24</w:t>
      </w:r>
    </w:p>'

$target.Range.InsertXML($newXml)

Write-Output ("Paragraphs now: " + $d.Paragraphs.Count)
